# Generate Report for Handoff
#
# Refresh the "Latest Handoff Datetime" (column D) for the files that were
# just (re-)handed off: 52ebef08-...-c7e5 (row 7), 646f79dc-...-1264 (row 11)
# and 8f751c9e-...-808f (row 13) on both the zh-cn and de-de status sheets.
# Rows 11/13 pick up the same refreshed timestamp as row 7 because their
# handoff ran as part of the same batch.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D7").Value  = "2016-03-09 02:37:37"
$zhcn.Range("D11").Value = "2016-03-09 02:37:37"
$zhcn.Range("D13").Value = "2016-03-09 02:37:37"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D7").Value  = "2016-03-09 02:37:48"
$dede.Range("D11").Value = "2016-03-09 02:37:48"
$dede.Range("D13").Value = "2016-03-09 02:37:48"
